# ---------------------------------------------------------------------------
# Applies the "minor changes, and updated Anleitung" edit:
#
#  Sheet "requests" (sheet1):
#   - Row 4 (Pulaski) and Row 8 (Culber) swap their entire B:AC schedule data
#     (the row label in column A stays with its own row).
#   - Q6 is cleared (was "N?").
#
#  Sheet "wards_min_presence" (sheet2):
#   - Row labelled "4_2" merges with "4_1" -> "4_2_4_1"; every filled cell in
#     that row becomes 2.
#   - Row labelled "12_2" merges with "6_3" -> "12_2_6_3"; every filled cell
#     in that row becomes 3.
#   - The old "4_1" row is removed (its values already lived on, relabelled,
#     on the row that used to say "4_3") and the old "6_3" row is removed.
#   - Net effect: rows 5 and 6 disappear, and the sheet shrinks from
#     A1:AC6 to A1:AC4.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsRequests = $wb.Worksheets.Item("requests")
$wsWards    = $wb.Worksheets.Item("wards_min_presence")

# ---------------------------------------------------------------------------
# 1) "requests": swap row 4 <-> row 8 (columns B..AC), keep column A as-is.
# ---------------------------------------------------------------------------

$dataCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$row4Values = @{}
$row8Values = @{}

foreach ($col in $dataCols) {
    $row4Values[$col] = $wsRequests.Range("${col}4").Value2
    $row8Values[$col] = $wsRequests.Range("${col}8").Value2
}

foreach ($col in $dataCols) {
    $v8 = $row8Values[$col]
    if ($v8 -eq $null -or $v8 -eq "") {
        $wsRequests.Range("${col}4").ClearContents()
    } else {
        $wsRequests.Range("${col}4").Value = $v8
    }

    $v4 = $row4Values[$col]
    if ($v4 -eq $null -or $v4 -eq "") {
        $wsRequests.Range("${col}8").ClearContents()
    } else {
        $wsRequests.Range("${col}8").Value = $v4
    }
}

# Clear Q6 (was "N?")
$wsRequests.Range("Q6").ClearContents()

# ---------------------------------------------------------------------------
# 2) "wards_min_presence": merge rows.
# ---------------------------------------------------------------------------

$presenceDataCols = @("B","E","F","G","H","I","L","M","N","O","P","S","T","U","V","W","Z","AA","AB","AC")

# Row 2 ("4_2") + row 3 ("4_1") -> row 2 becomes "4_2_4_1" with every filled
# cell set to 2.
$wsWards.Range("A2").Value = "4_2_4_1"
foreach ($col in $presenceDataCols) {
    $wsWards.Range("${col}2").Value = 2
}

# Row 4 ("12_2") + row 6 ("6_3") -> row 4 becomes "12_2_6_3" with every
# filled cell set to 3.
$wsWards.Range("A4").Value = "12_2_6_3"
foreach ($col in $presenceDataCols) {
    $wsWards.Range("${col}4").Value = 3
}

# Remove the now-redundant rows. Row 3 keeps its old ("4_1") data but is
# relabelled "4_3" by the row that follows shifting up, so delete rows 6 and
# 5 (highest index first so row numbers of not-yet-deleted rows don't shift).
$wsWards.Rows.Item(6).Delete()
$wsWards.Rows.Item(5).Delete()

$wsWards.Range("A3").Value = "4_3"
